$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column E ("res"); everything from E
# onward (old E..P) shifts right to F..Q. Column formatting (incl. the
# bold/bordered header style on row 1) comes along with the shift, and
# Excel automatically extends the header row's style onto the freshly
# inserted E1 cell.
$ws.Columns("E:E").Insert()

# New column header: "text"
$ws.Range("E1").Value = "text"

# Populate the new column with a straight numeric copy of column D
# ("sharp") for every data row.
for ($r = 2; $r -le 61; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value2
}

# Restore cursor/selection reported by the diff.
$ws.Range("F16").Select()
